$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.536.06"
$ws.Range("E2").Value = "  +3.72%  "
$ws.Range("D3").Value = "2.439.91"
$ws.Range("E3").Value = "  +2.39%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "573.77"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.98%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.71"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.82%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.08%  "
$ws.Range("E8").Value = "  +1.57%  "
$ws.Range("D9").Value = "2.439.65"
$ws.Range("E9").Value = "  +2.35%  "
$ws.Range("E10").Value = "  +4.54%  "
$ws.Range("E11").Value = "  +0.70%  "
$ws.Range("E12").Value = "  +2.05%  "
$ws.Range("E13").Value = "  +3.76%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.42"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +5.98%  "
$ws.Range("E15").Value = "  +5.45%  "
$ws.Range("D16").Value = "2.883.77"
$ws.Range("E16").Value = "  +2.50%  "
$ws.Range("D17").Value = "62.469.90"
$ws.Range("E17").Value = "  +3.88%  "
$ws.Range("D18").Value = "2.438.06"
$ws.Range("E18").Value = "  +2.08%  "
$ws.Range("E19").Value = "  -3.03%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.91"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.44%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "327.44"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.42%  "
$ws.Range("E23").Value = "  +11.75%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.00"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.21%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "65.57"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.01%  "
$ws.Range("D26").Value = "0.0₆0610"
$ws.Range("E26").Value = "  +109.27%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "623.22"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +10.59%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.11"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +10.57%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.53"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +5.11%  "
$ws.Range("D31").Value = "2.561.13"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.17"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.16%  "
$ws.Range("E33").Value = "  +6.66%  "
$ws.Range("E34").Value = "  +3.55%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.84"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.38%  "
$ws.Range("E36").Value = "  +2.52%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.999"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.08%  "
$ws.Range("E38").Value = "  +3.50%  "
$ws.Range("E39").Value = "  +1.48%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "151.88"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.07%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.39"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +5.95%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "18.59"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.16%  "
$ws.Range("E43").Value = "  +13.33%  "
$ws.Range("E44").Value = "  +4.30%  "
$ws.Range("E45").Value = "  +0.01%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "14.72"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +25.71%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "144.32"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.69%  "
$ws.Range("E48").Value = "  +1.54%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "20.50"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +6.72%  "
$ws.Range("E50").Value = "  +1.60%  "
$ws.Range("E51").Value = "  +2.22%  "
